$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 442, shifting the
# existing rows 442..538 down to 444..540 (matches Excel's normal
# "insert rows" shift-down behaviour).
$ws.Rows("442:443").Insert()

# Row 442 (new price quote for Ají / Inferno / Primera, Región de Arica y
# Parinacota, $/caja 10 kilos)
$ws.Cells.Item(442, 1).Value = 8
$ws.Cells.Item(442, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(442, 3).Value = "Coquimbo"
$ws.Cells.Item(442, 4).Value = 45173
$ws.Cells.Item(442, 5).Value = 4
$ws.Cells.Item(442, 6).Value = 100112021
$ws.Cells.Item(442, 7).Value = "Ají"
$ws.Cells.Item(442, 8).Value = "Inferno"
$ws.Cells.Item(442, 9).Value = "Primera"
$ws.Cells.Item(442, 10).Value = 500
$ws.Cells.Item(442, 11).Value = 20000
$ws.Cells.Item(442, 12).Value = 21000
$ws.Cells.Item(442, 13).Value = 20500
$ws.Cells.Item(442, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(442, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(442, 16).Value = 2050
$ws.Cells.Item(442, 17).Value = 10
$ws.Cells.Item(442, 18).Value = "Hortaliza"

# Row 443 (new price quote for Ají / Inferno / Segunda, Región de Arica y
# Parinacota, $/caja 10 kilos)
$ws.Cells.Item(443, 1).Value = 8
$ws.Cells.Item(443, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(443, 3).Value = "Coquimbo"
$ws.Cells.Item(443, 4).Value = 45173
$ws.Cells.Item(443, 5).Value = 4
$ws.Cells.Item(443, 6).Value = 100112021
$ws.Cells.Item(443, 7).Value = "Ají"
$ws.Cells.Item(443, 8).Value = "Inferno"
$ws.Cells.Item(443, 9).Value = "Segunda"
$ws.Cells.Item(443, 10).Value = 320
$ws.Cells.Item(443, 11).Value = 13000
$ws.Cells.Item(443, 12).Value = 14000
$ws.Cells.Item(443, 13).Value = 13500
$ws.Cells.Item(443, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(443, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(443, 16).Value = 1350
$ws.Cells.Item(443, 17).Value = 10
$ws.Cells.Item(443, 18).Value = "Hortaliza"
